# Apply the edit described by the diff:
#  - Fill in row 3 (B3, C3, E3, F3, G3) mirroring row 2's pattern
#  - Add a new row 4 with C4 = 0
#  - The dimension will grow from A1:G3 to A1:G4 automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = $false
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 19.170000000000002
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = $false

$ws.Range("C4").Value = 0
